$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number must be pre-formatted as
# Text so Excel keeps storing/ echoing them as the literal inline string from
# the source data (matching the original t="inlineStr" cells) instead of
# silently converting them to a floating point number on entry.
$numericLookingCells = @('D5', 'D6', 'D10', 'D16', 'D19', 'D21', 'D22', 'D26', 'D27', 'D28', 'D32', 'D35', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D45', 'D46', 'D49', 'D50')
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '61.137.91'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').Value = '3.386.43'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '574.26'
$ws.Range('D6').Value = '137.21'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D8').Value = '3.385.51'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  -1.12%  '
$ws.Range('D10').Value = '7.63'
$ws.Range('E10').Value = '  +1.91%  '
$ws.Range('E11').Value = '  -2.96%  '
$ws.Range('E12').Value = '  -2.39%  '
$ws.Range('D13').Value = '3.954.36'
$ws.Range('E13').Value = '  -0.44%  '
$ws.Range('E14').Value = '  +0.69%  '
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('D16').Value = '25.69'
$ws.Range('E16').Value = '  +1.33%  '
$ws.Range('D17').Value = '3.385.22'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').Value = '61.276.93'
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('D19').Value = '13.82'
$ws.Range('E19').Value = '  -2.20%  '
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('D21').Value = '9.33'
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('D22').Value = '376.48'
$ws.Range('E22').Value = '  -0.78%  '
$ws.Range('D23').Value = '3.519.53'
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('E24').Value = '  -2.19%  '
$ws.Range('E25').Value = '  +0.21%  '
$ws.Range('D26').Value = '0.0000125'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('D27').Value = '70.99'
$ws.Range('E27').Value = '  -0.38%  '
$ws.Range('D28').Value = '0.178'
$ws.Range('E28').Value = '  +11.53%  '
$ws.Range('E29').Value = '  -3.74%  '
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('E31').Value = '  -2.18%  '
$ws.Range('D32').Value = '8.07'
$ws.Range('E32').Value = '  -1.65%  '
$ws.Range('E33').Value = '  -1.66%  '
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').Value = '23.43'
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('E36').Value = '  -4.17%  '
$ws.Range('E37').Value = '  -1.57%  '
$ws.Range('D38').Value = '6.83'
$ws.Range('E38').Value = '  -0.64%  '
$ws.Range('D39').Value = '164.73'
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('D40').Value = '0.0762'
$ws.Range('E40').Value = '  -3.21%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '25.52'
$ws.Range('E41').Value = '  +2.29%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').Value = '0.775'
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('E44').Value = '  -1.43%  '
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').Value = '1.19'
$ws.Range('E45').Value = '  -3.55%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '4.34'
$ws.Range('E46').Value = '  -1.78%  '
$ws.Range('D47').Value = '2.558.27'
$ws.Range('E47').Value = '  +9.01%  '
$ws.Range('E48').Value = '  -1.17%  '
$ws.Range('D49').Value = '22.93'
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('D50').Value = '2.44'
$ws.Range('E50').Value = '  +4.15%  '
$ws.Range('E51').Value = '  -1.33%  '

# Restore the default (General) style on those cells so no stray number
# formatting lingers on the cell itself once the text is safely stored.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
